$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A to text format first so the date-like string "2025-10-11"
# is stored as a literal string (matching the inlineStr in the diff) rather
# than being auto-converted by Excel into a date serial number.
$ws.Range("A57").NumberFormat = "@"
$ws.Range("A57").Value = "2025-10-11"
# Reset the cell style back to the default/normal style so no extra
# formatting (e.g. the text number format) lingers on the new cell,
# matching the unstyled cells used by the rest of the data rows.
$ws.Range("A57").Style = "Normal"

$ws.Range("B57").Value = 54.31000137329102
$ws.Range("C57").Value = 678.9500122070312
$ws.Range("D57").Value = 348.2999877929688
